$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.860.50'
$ws.Range("E2").Value = '  +6.87%  '
$ws.Range("D3").Value = '2.622.63'
$ws.Range("E3").Value = '  +8.63%  '
$ws.Range("E4").Value = '  -0.39%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '507.58'
$ws.Range("E5").Value = '  +3.74%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '157.55'
$ws.Range("E6").Value = '  +1.59%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.995'
$ws.Range("E7").Value = '  -0.25%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.591'
$ws.Range("E8").Value = '  -4.15%  '
$ws.Range("D9").Value = '2.664.88'
$ws.Range("E9").Value = '  +9.52%  '
$ws.Range("E10").Value = '  +1.65%  '
$ws.Range("E11").Value = '  +4.82%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.344'
$ws.Range("E12").Value = '  +3.37%  '
$ws.Range("E13").Value = '  +0.88%  '
$ws.Range("D14").Value = '3.089.56'
$ws.Range("E14").Value = '  +8.75%  '
$ws.Range("D15").Value = '60.797.70'
$ws.Range("E15").Value = '  +6.49%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.76'
$ws.Range("E16").Value = '  +5.48%  '
$ws.Range("E17").Value = '  +5.96%  '
$ws.Range("D18").Value = '2.658.83'
$ws.Range("E18").Value = '  +9.65%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.79'
$ws.Range("E19").Value = '  +1.67%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '346.21'
$ws.Range("E20").Value = '  +6.57%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.52'
$ws.Range("E21").Value = '  +4.99%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.19'
$ws.Range("E22").Value = '  +4.43%  '
$ws.Range("E23").Value = '  -0.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '60.57'
$ws.Range("E24").Value = '  +4.92%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.424'
$ws.Range("E25").Value = '  +4.83%  '
$ws.Range("E26").Value = '  +3.92%  '
$ws.Range("E27").Value = '  -0.13%  '
$ws.Range("D28").Value = '0.0₃0864'
$ws.Range("E28").Value = '  +10.61%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.61'
$ws.Range("E29").Value = '  +4.80%  '
$ws.Range("E30").Value = '  -0.27%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '19.54'
$ws.Range("E31").Value = '  +4.79%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '156.05'
$ws.Range("E32").Value = '  +3.29%  '
$ws.Range("E33").Value = '  +3.00%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.79'
$ws.Range("E34").Value = '  +9.43%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.05'
$ws.Range("E35").Value = '  +7.10%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.21'
$ws.Range("E36").Value = '  +5.21%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '310.18'
$ws.Range("E37").Value = '  +10.34%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.857'
$ws.Range("E38").Value = '  +3.29%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.49'
$ws.Range("E39").Value = '  +8.78%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.839'
$ws.Range("E40").Value = '  +29.64%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.77'
$ws.Range("E41").Value = '  +6.76%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '35.40'
$ws.Range("E42").Value = '  +4.05%  '
$ws.Range("E43").Value = '  +6.62%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0576'
$ws.Range("E44").Value = '  +8.39%  '
$ws.Range("E45").Value = '  -2.17%  '
$ws.Range("E46").Value = '  +14.93%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.996'
$ws.Range("E47").Value = '  +0.16%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.92'
$ws.Range("E48").Value = '  +9.04%  '
$ws.Range("E49").Value = '  +4.44%  '
$ws.Range("D50").Value = '2.054.04'
$ws.Range("E50").Value = '  +8.14%  '
$ws.Range("E51").Value = '  +0.22%  '
